$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1350.6666
$ws.Range("I19").Value = 1173
$ws.Range("J19").Value = 1706
$ws.Range("K19").Value = 1173
$ws.Range("L19").Value = 1706
$ws.Range("M19").Value = -998
$ws.Range("N19").Value = -2056
$ws.Range("H96").Value = 4464825
$ws.Range("I96").Value = 7936734
$ws.Range("K96").Value = 23810202
$ws.Range("M96").Value = -23808829
$ws.Range("H98").Value = 6347.25
$ws.Range("I98").Value = 3129.5833
$ws.Range("J98").Value = 16000.25
$ws.Range("K98").Value = 3129.5833
$ws.Range("L98").Value = 16000.25
$ws.Range("M98").Value = -1631.5833
$ws.Range("N98").Value = -18996.25
$ws.Range("H122").Value = 6347.25
$ws.Range("I122").Value = 3129.5833
$ws.Range("J122").Value = 16000.25
$ws.Range("K122").Value = 9388.749899999999
$ws.Range("L122").Value = 48000.75
$ws.Range("M122").Value = -6938.749899999999
$ws.Range("N122").Value = -52900.75
$ws.Range("H125").Value = 15910518
$ws.Range("I125").Value = 47722388
$ws.Range("J125").Value = 4582.6665
$ws.Range("K125").Value = 429501492
$ws.Range("L125").Value = 41243.9985
$ws.Range("M125").Value = -429499032
$ws.Range("N125").Value = -46163.9985
$ws.Range("H132").Value = 38613.58
$ws.Range("I132").Value = 38613.58
$ws.Range("K132").Value = 115840.74
$ws.Range("M132").Value = -113310.74

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 306005.6
$ws.Range("J34").Value = 257507
$ws.Range("L34").Value = 257507
$ws.Range("N34").Value = -258049
$ws.Range("H61").Value = 13141.333
$ws.Range("I61").Value = 1591.8
$ws.Range("J61").Value = 27578.25
$ws.Range("K61").Value = 1591.8
$ws.Range("L61").Value = 27578.25
$ws.Range("M61").Value = -1379.8
$ws.Range("N61").Value = -28002.25
$ws.Range("H110").Value = 547.7857
$ws.Range("I110").Value = 513
$ws.Range("K110").Value = 513
$ws.Range("M110").Value = 1532
$ws.Range("H132").Value = 1036.7059
$ws.Range("I132").Value = 1036.7059
$ws.Range("K132").Value = 3110.1177
$ws.Range("M132").Value = -580.1176999999998
$ws.Range("H136").Value = 13141.333
$ws.Range("I136").Value = 1591.8
$ws.Range("J136").Value = 27578.25
$ws.Range("K136").Value = 4775.4
$ws.Range("L136").Value = 82734.75
$ws.Range("M136").Value = -2225.4
$ws.Range("N136").Value = -87834.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 23968.8
$ws.Range("I75").Value = 16614.666
$ws.Range("J75").Value = 35000
$ws.Range("K75").Value = 16614.666
$ws.Range("L75").Value = 35000
$ws.Range("M75").Value = -15678.666
$ws.Range("N75").Value = -36872
$ws.Range("H78").Value = 23968.8
$ws.Range("I78").Value = 16614.666
$ws.Range("J78").Value = 35000
$ws.Range("K78").Value = 49843.99800000001
$ws.Range("L78").Value = 105000
$ws.Range("M78").Value = -45163.99800000001
$ws.Range("N78").Value = -114360
$ws.Range("H99").Value = 1378.1875
$ws.Range("I99").Value = 1254.5
$ws.Range("K99").Value = 1254.5
$ws.Range("M99").Value = 243.5
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("H107").Value = 4595.8423
$ws.Range("I107").Value = 4954.75
$ws.Range("K107").Value = 4954.75
$ws.Range("M107").Value = -3034.75
$ws.Range("H123").Value = 40454.547
$ws.Range("J123").Value = 37857.145
$ws.Range("L123").Value = 37857.145
$ws.Range("N123").Value = -47657.145
$ws.Range("H134").Value = 9729.156000000001
$ws.Range("I134").Value = 10494.296
$ws.Range("J134").Value = 5597.4
$ws.Range("K134").Value = 31482.888
$ws.Range("L134").Value = 16792.2
$ws.Range("M134").Value = -28947.888
$ws.Range("N134").Value = -21862.2
$ws.Range("N101").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2141.4
$ws.Range("I16").Value = 1426.875
$ws.Range("K16").Value = 1426.875
$ws.Range("M16").Value = -1139.875
$ws.Range("H96").Value = 5500
$ws.Range("J96").Value = 5500
$ws.Range("L96").Value = 5500
$ws.Range("N96").Value = -10992
$ws.Range("H99").Value = 6497.3335
$ws.Range("I99").Value = 7857
$ws.Range("K99").Value = 7857
$ws.Range("M99").Value = -6359
$ws.Range("H113").Value = 2141.4
$ws.Range("I113").Value = 1426.875
$ws.Range("K113").Value = 1426.875
$ws.Range("M113").Value = 743.125
$ws.Range("H126").Value = 6497.3335
$ws.Range("I126").Value = 7857
$ws.Range("K126").Value = 23571
$ws.Range("M126").Value = -21101

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2317.904
$ws.Range("I131").Value = 4114.9
$ws.Range("J131").Value = 2032.6666
$ws.Range("K131").Value = 12344.7
$ws.Range("L131").Value = 6097.9998
$ws.Range("M131").Value = -7304.699999999999
$ws.Range("N131").Value = -16177.9998

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 26600
$ws.Range("J74").Value = 26600
$ws.Range("L74").Value = 26600
$ws.Range("N74").Value = -28472
$ws.Range("H77").Value = 26600
$ws.Range("J77").Value = 26600
$ws.Range("L77").Value = 79800
$ws.Range("N77").Value = -89160
$ws.Range("H80").Value = 4090.4517
$ws.Range("I80").Value = 2631.318
$ws.Range("J80").Value = 7657.222
$ws.Range("K80").Value = 2631.318
$ws.Range("L80").Value = 7657.222
$ws.Range("M80").Value = -1633.318
$ws.Range("N80").Value = -9653.222
$ws.Range("H83").Value = 4090.4517
$ws.Range("I83").Value = 2631.318
$ws.Range("J83").Value = 7657.222
$ws.Range("K83").Value = 13156.59
$ws.Range("L83").Value = 38286.11
$ws.Range("M83").Value = -8164.59
$ws.Range("N83").Value = -48270.11
$ws.Range("H113").Value = 3661
$ws.Range("I113").Value = 1652.5
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 1652.5
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 517.5
$ws.Range("N113").Value = -9340
$ws.Range("H123").Value = 41562.312
$ws.Range("J123").Value = 41562.312
$ws.Range("L123").Value = 41562.312
$ws.Range("N123").Value = -46462.312
$ws.Range("H126").Value = 2351.6
$ws.Range("I126").Value = 2142.5715
$ws.Range("K126").Value = 6427.7145
$ws.Range("M126").Value = -3957.7145
$ws.Range("H132").Value = 2780.2856
$ws.Range("I132").Value = 2437.3572
$ws.Range("K132").Value = 7312.071599999999
$ws.Range("M132").Value = -4782.071599999999
$ws.Range("H139").Value = 115000
$ws.Range("J139").Value = 115000
$ws.Range("L139").Value = 115000
$ws.Range("N139").Value = -125280

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 13718.667
$ws.Range("J42").Value = 11671.125
$ws.Range("L42").Value = 11671.125
$ws.Range("N42").Value = -12797.125
$ws.Range("H43").Value = 16557.182
$ws.Range("J43").Value = 15001.9
$ws.Range("L43").Value = 15001.9
$ws.Range("N43").Value = -15387.9
$ws.Range("H46").Value = 4093.2307
$ws.Range("I46").Value = 5000
$ws.Range("K46").Value = 5000
$ws.Range("M46").Value = -4812
$ws.Range("H49").Value = 13718.667
$ws.Range("J49").Value = 11671.125
$ws.Range("L49").Value = 11671.125
$ws.Range("N49").Value = -11965.125
$ws.Range("H55").Value = 3147.25
$ws.Range("J55").Value = 4063
$ws.Range("L55").Value = 4063
$ws.Range("N55").Value = -4409
$ws.Range("H100").Value = 2489.9583
$ws.Range("I100").Value = 2414.6667
$ws.Range("K100").Value = 2414.6667
$ws.Range("M100").Value = -1873.6667
$ws.Range("H122").Value = 3925
$ws.Range("I122").Value = 3900
$ws.Range("J122").Value = 3933.3333
$ws.Range("K122").Value = 11700
$ws.Range("L122").Value = 11799.9999
$ws.Range("M122").Value = -9250
$ws.Range("N122").Value = -16699.9999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 118499.5
$ws.Range("J46").Value = 118499.5
$ws.Range("L46").Value = 118499.5
$ws.Range("N46").Value = -118961.5
$ws.Range("H122").Value = 92186.06
$ws.Range("I122").Value = 104727.29
$ws.Range("J122").Value = 4397.5
$ws.Range("K122").Value = 314181.87
$ws.Range("L122").Value = 13192.5
$ws.Range("M122").Value = -311731.87
$ws.Range("N122").Value = -18092.5
$ws.Range("H134").Value = 118499.5
$ws.Range("J134").Value = 118499.5
$ws.Range("L134").Value = 355498.5
$ws.Range("N134").Value = -360568.5
$ws.Range("H136").Value = 24660.74
$ws.Range("I136").Value = 31512.45
$ws.Range("J136").Value = 5084.4287
$ws.Range("K136").Value = 94537.35000000001
$ws.Range("L136").Value = 15253.2861
$ws.Range("M136").Value = -91987.35000000001
$ws.Range("N136").Value = -20353.2861
